# Update the CFU summary table (medians per vendor, per day) with newly
# computed median values (plotting medians instead of a single overall
# median), per Nick's feedback.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2" = 2000
    "G2" = 1000
    "H2" = 195000
    "I2" = 0
    "J2" = 100000
    "K2" = 17500

    "F3" = 0
    "G3" = 0
    "H3" = 92000
    "I3" = 0
    "J3" = 340000

    "F4" = 0
    "G4" = 0
    "H4" = 65000
    "I4" = 0
    "J4" = 12000
    "K4" = 1400

    "F5" = 60000
    "G5" = 180000
    "H5" = 4800000
    "I5" = 120000
    "J5" = 1305000
    "K5" = 1770000

    "F6" = 400000
    "G6" = 1235000
    "H6" = 22100000
    "I6" = 960000
    "J6" = 400000
    "K6" = 9670000

    "J7" = 2000
    "K7" = 0

    "F8" = 19800000
    "G8" = 6300000
    "I8" = 420000
    "K8" = 52200000

    "F9" = 59500000
    "G9" = 53000000
    "H9" = 126000000
    "I9" = 33000000
    "J9" = 38500000
    "K9" = 22000000

    "J10" = 0
    "K10" = 0
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
